$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new column M that duplicates column A (Phenotype names) for rows 1-22,
# so the phenotype label stays visible alongside the new BLUP-range columns.
$ws.Range("A1:A22").Copy()
$ws.Range("M1").PasteSpecial(-4163)
$excel.CutCopyMode = 0

# Turn the per-row I formulas (=H*D) into one shared formula spanning I2:I22.
$ws.Range("I2:I22").Formula = "=H2*D2"

# Turn the contiguous per-row J formulas (=I/(G-F)) for rows 6-14 into a shared formula.
$ws.Range("J6:J14").Formula = "=I6/(G6-F6)"

# Update the view so column I is the left-most visible column and the single
# active cell is L14 (previously K2:K14 was selected).
$ws.Range("L14").Select()
